$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Globo", "RJ TV 2", "Defesa Civil", "2025-04-02T18:44", "Neutro", "Retirada dos escombros de prédio em Campos. Imóvel histórico no Centro precisou ser demolido depois que parte do teto cedeu. Defesa Civil interditou o prédio na semana passada. Repórter *ao vivo*"),
    @("Globo", "RJ TV 2", "Transporte", "2025-04-02T18:54", "Negativo", "Por whatsApp, reclamação de falta de ônibus no Parque Imperial"),
    @("Globo", "RJ TV 2", "Educação", "2025-04-02T18:55", "Negativo", "Por whatsApp, moradora pede atenção à Creche Wilson Amaro Freitas, no Jockey. Ano letivo começou há dois meses e a poda na área interna ainda não foi realizada. Ofício já foi enviado. "),
    @("Globo", "RJ TV 2", "CCZ", "2025-04-02T18:57", "Positivo", "Animais abandonados em Campos. CCZ faz alerta para perigos de animais soltos em vias públicas. Repórter *ao vivo*. Depoimento do diretor do CCZ, Carlos Morales. ")
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
